# Official_data/Variables.xlsx - "Add files via upload"
#
# The sheet gains a "Description" header in C1 (column C already held the
# short variable descriptions; it previously had no header text there),
# and the long-form description in C2 (the Life Expectancy row) has its
# wrap-text turned off so it behaves like a single-line cell instead of a
# wrapped paragraph. The active selection is also moved from C14 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the (already-present) description column.
$ws.Range("C1").Value = "Description"

# Turn off word-wrap for C2 only (keeps its top-vertical alignment).
$ws.Range("C2").WrapText = $false

# Move the selection / active cell to C2.
$ws.Range("C2").Select() | Out-Null
